# Generate Report for Handoff
# Updates the row for "a788d847-8a9b-4446-992f-1622740ddb20.md" across the
# Overview / zh-cn / de-de sheets: status flips from "Handed back: in sync
# with en-US" to "Ready for handoff", the xliff generation timestamps move
# forward, and an Error Detail message is recorded for zh-cn/de-de noting
# the handback file version is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9bff73dc443b49c71c25f7da366e50e13dc897cb/e2e/a788d847-8a9b-4446-992f-1622740ddb20.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69a379728fd3d5a9b409129bd0dd53b6cb373c7d/e2e/a788d847-8a9b-4446-992f-1622740ddb20.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-20 18:57:17"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-20 18:57:13"
$wsZhCn.Range("P3").Value = $errorDetail
# Native width 40 (ColumnWidth read-back/round-trip value for width=40)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-20 18:57:17"
$wsDeDe.Range("P3").Value = $errorDetail
# Native width 40 (ColumnWidth read-back/round-trip value for width=40)
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
